$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows before row 6; this shifts the existing rows 6-11 down to 8-13
$ws.Range("A6:R7").EntireRow.Insert()

# Match the number format/style of the date column (D) used by the other data rows
$ws.Range("D6:D7").NumberFormat = $ws.Range("D8").NumberFormat

# Fill out row 6 with the new data
$ws.Cells.Item(6, 1).Value = 7
$ws.Cells.Item(6, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(6, 3).Value = "Ñuble"
$ws.Cells.Item(6, 4).Value = 44791
$ws.Cells.Item(6, 5).Value = 16
$ws.Cells.Item(6, 6).Value = 100112044
$ws.Cells.Item(6, 7).Value = "Perejil"
$ws.Cells.Item(6, 8).Value = "Sin especificar"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 240
$ws.Cells.Item(6, 11).Value = 750
$ws.Cells.Item(6, 12).Value = 800
$ws.Cells.Item(6, 13).Value = 775
$ws.Cells.Item(6, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(6, 15).Value = "Región del Maule"
$ws.Cells.Item(6, 16).Value = 775
$ws.Cells.Item(6, 17).Value = 1
$ws.Cells.Item(6, 18).Value = "Hortaliza"

# Fill out row 7 with the new data
$ws.Cells.Item(7, 1).Value = 7
$ws.Cells.Item(7, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(7, 3).Value = "Ñuble"
$ws.Cells.Item(7, 4).Value = 44791
$ws.Cells.Item(7, 5).Value = 16
$ws.Cells.Item(7, 6).Value = 100112044
$ws.Cells.Item(7, 7).Value = "Perejil"
$ws.Cells.Item(7, 8).Value = "Sin especificar"
$ws.Cells.Item(7, 9).Value = "Segunda"
$ws.Cells.Item(7, 10).Value = 250
$ws.Cells.Item(7, 11).Value = 650
$ws.Cells.Item(7, 12).Value = 650
$ws.Cells.Item(7, 13).Value = 650
$ws.Cells.Item(7, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(7, 15).Value = "Región del Maule"
$ws.Cells.Item(7, 16).Value = 650
$ws.Cells.Item(7, 17).Value = 1
$ws.Cells.Item(7, 18).Value = "Hortaliza"
